# Commit: "add estilos boostrap, imagen arquitectura"
#
# Structural change visible in the OOXML diff:
#   - A new worksheet "hola" is appended after "PINI".
#   - "papa" (sheet1) is trimmed down to just its header row
#     (the "Auricular Inalambrico" / "5" data row is removed).
#   - "PINI" (sheet2) gains a data row: "MASO" / "4"
#     (previously it only had the header row).
#   - The new "hola" sheet gets the same header-only layout
#     ("Producto" / "Cantidad").

$wb = $excel.ActiveWorkbook

$papa = $wb.Worksheets.Item("papa")
$pini = $wb.Worksheets.Item("PINI")

# --- "papa": drop the data row, keep only the header row ---
$papa.Range("A2:B2").ClearContents()

# --- "PINI": add the new data row (kept as text, matching the sibling
#     "5" cell already stored as text elsewhere in this workbook) ---
$pini.Range("A2").Value = "MASO"
$pini.Range("B2").NumberFormat = "@"
$pini.Range("B2").Value = "4"

# --- add the new "hola" worksheet right after "PINI" ---
$hola = $wb.Worksheets.Add($null, $pini)
$hola.Name = "hola"
$hola.Range("A1").Value = "Producto"
$hola.Range("B1").Value = "Cantidad"
